$d = $word.ActiveDocument

# --- Heading 1 & 2: turn off "contextual spacing" (w:contextualSpacing) ---
foreach ($name in @("Heading1", "Heading2")) {
    $s = $d.Styles($name)
    $s.NoSpaceBetweenParagraphsOfSameStyle = $false
}

# --- Heading 3, 4, 5: rebase on Normal, turn off contextual spacing, and
#     bump the run font size to 12pt (w:sz 24 half-points) ---
foreach ($name in @("Heading3", "Heading4", "Heading5")) {
    $s = $d.Styles($name)
    $s.BaseStyle = $d.Styles("Normal")
    $s.NoSpaceBetweenParagraphsOfSameStyle = $false
    $s.Font.Size = 12
}
